$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content changes ---
# D30: shorten the ETC identification-code description text
$ws.Range("D30").Value = "Ebola Care Facility's identification code"

# A28:A32 - fill in the "SOURCE" column for the Organization attribute rows
# (3W Dataset on the Organizations Involved in the Response to the Ebola Crisis)
$orgSource = "3W Dataset on the Organizations Involved in the Response to the Ebola Crisis"
$ws.Range("A28").Value = $orgSource
$ws.Range("A29").Value = $orgSource
$ws.Range("A30").Value = $orgSource
$ws.Range("A31").Value = $orgSource
$ws.Range("A32").Value = $orgSource

# --- Column width changes ---
# Column A widened (and best-fit) to accommodate the longer attribute names
$ws.Columns.Item(1).ColumnWidth = 62.8333333333333
# Columns B and D revert to the (default) width, since they no longer need a custom width
$ws.Columns.Item(2).ColumnWidth = 9.67
$ws.Columns.Item(4).ColumnWidth = 9.67

# --- Sheet view changes ---
# Scroll so row 4 is the top visible row (topLeftCell A4), then select A13
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("A13").Select()

# --- Workbook window position ---
$wb.Windows.Item(1).Left = 10820
